$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 366
$ws.Range("I12").Value = 392
$ws.Range("J12").Value = 305.33334
$ws.Range("K12").Value = 392
$ws.Range("L12").Value = 305.33334
$ws.Range("M12").Value = -222
$ws.Range("N12").Value = -645.33334
$ws.Range("H17").Value = 592461.7
$ws.Range("J17").Value = 626990.5600000001
$ws.Range("L17").Value = 1880971.68
$ws.Range("N17").Value = -1881307.68
$ws.Range("H41").Value = 111596.445
$ws.Range("I41").Value = 531.25
$ws.Range("J41").Value = 200448.6
$ws.Range("K41").Value = 531.25
$ws.Range("L41").Value = 200448.6
$ws.Range("M41").Value = -91.25
$ws.Range("N41").Value = -201328.6
$ws.Range("H107").Value = 4533.8823
$ws.Range("I107").Value = 519.5333000000001
$ws.Range("K107").Value = 519.5333000000001
$ws.Range("M107").Value = 1400.4667
$ws.Range("H113").Value = 125003750
$ws.Range("I113").Value = 333336000
$ws.Range("J113").Value = 4400
$ws.Range("K113").Value = 333336000
$ws.Range("L113").Value = 4400
$ws.Range("M113").Value = -333332746
$ws.Range("N113").Value = -10908
$ws.Range("H125").Value = 3505.9092
$ws.Range("I125").Value = 3646.1667
$ws.Range("J125").Value = 3337.6
$ws.Range("K125").Value = 32815.5003
$ws.Range("L125").Value = 30038.4
$ws.Range("M125").Value = -30355.5003
$ws.Range("N125").Value = -34958.39999999999
$ws.Range("H137").Value = 2373.625
$ws.Range("I137").Value = 1738.52
$ws.Range("K137").Value = 5215.559999999999
$ws.Range("M137").Value = -2665.559999999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13752.967
$ws.Range("I32").Value = 8199.933999999999
$ws.Range("J32").Value = 19306
$ws.Range("K32").Value = 8199.933999999999
$ws.Range("L32").Value = 19306
$ws.Range("M32").Value = -7912.933999999999
$ws.Range("N32").Value = -19880
$ws.Range("H61").Value = 1528.5161
$ws.Range("I61").Value = 1329.4667
$ws.Range("K61").Value = 1329.4667
$ws.Range("M61").Value = -1117.4667
$ws.Range("H97").Value = 627.2549
$ws.Range("I97").Value = 604.35
$ws.Range("J97").Value = 710.5454999999999
$ws.Range("K97").Value = 604.35
$ws.Range("L97").Value = 710.5454999999999
$ws.Range("M97").Value = -108.35
$ws.Range("N97").Value = -1702.5455
$ws.Range("H101").Value = 8000
$ws.Range("J101").Value = 8000
$ws.Range("L101").Value = 8000
$ws.Range("N101").Value = -14490
$ws.Range("H122").Value = 2361.2
$ws.Range("I122").Value = 1570
$ws.Range("K122").Value = 4710
$ws.Range("M122").Value = -2260
$ws.Range("H132").Value = 2790.1304
$ws.Range("I132").Value = 2249.1333
$ws.Range("J132").Value = 3804.5
$ws.Range("K132").Value = 6747.3999
$ws.Range("L132").Value = 11413.5
$ws.Range("M132").Value = -4217.3999
$ws.Range("N132").Value = -16473.5
$ws.Range("H136").Value = 1528.5161
$ws.Range("I136").Value = 1329.4667
$ws.Range("K136").Value = 3988.4001
$ws.Range("M136").Value = -1438.4001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1446.238
$ws.Range("I105").Value = 1114.0667
$ws.Range("J105").Value = 2276.6667
$ws.Range("K105").Value = 1114.0667
$ws.Range("L105").Value = 2276.6667
$ws.Range("M105").Value = 632.9332999999999
$ws.Range("N105").Value = -5770.6667
$ws.Range("H126").Value = 65099
$ws.Range("J126").Value = 65099
$ws.Range("L126").Value = 65099
$ws.Range("N126").Value = -74979
$ws.Range("H134").Value = 1544.2693
$ws.Range("I134").Value = 1408.04
$ws.Range("K134").Value = 4224.12
$ws.Range("M134").Value = -1689.12

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 30761
$ws.Range("J28").Value = 30761
$ws.Range("L28").Value = 30761
$ws.Range("N28").Value = -31251
$ws.Range("H43").Value = 19599.75
$ws.Range("J43").Value = 19599.75
$ws.Range("L43").Value = 19599.75
$ws.Range("N43").Value = -19967.75
$ws.Range("H94").Value = 1749.6666
$ws.Range("I94").Value = 1511
$ws.Range("J94").Value = 1797.4
$ws.Range("K94").Value = 1511
$ws.Range("L94").Value = 1797.4
$ws.Range("M94").Value = -1060
$ws.Range("N94").Value = -2699.4
$ws.Range("H101").Value = 19599.75
$ws.Range("J101").Value = 19599.75
$ws.Range("L101").Value = 19599.75
$ws.Range("N101").Value = -26089.75
$ws.Range("H132").Value = 3116.4565
$ws.Range("I132").Value = 2996.5952
$ws.Range("K132").Value = 8989.785600000001
$ws.Range("M132").Value = -6459.785600000001
$ws.Range("H134").Value = 2301.0833
$ws.Range("I134").Value = 2229.4348
$ws.Range("K134").Value = 6688.3044
$ws.Range("M134").Value = -4153.3044
$ws.Range("H141").Value = 229165.73
$ws.Range("J141").Value = 229165.73
$ws.Range("L141").Value = 229165.73
$ws.Range("N141").Value = -239525.73

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 25000
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H88").Value = 10333.333
$ws.Range("J88").Value = 10333.333
$ws.Range("L88").Value = 30999.999
$ws.Range("N88").Value = -31855.999
$ws.Range("H91").Value = 10333.333
$ws.Range("J91").Value = 10333.333
$ws.Range("L91").Value = 30999.999
$ws.Range("N91").Value = -33963.999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2112.25
$ws.Range("I102").Value = 1676.6923
$ws.Range("J102").Value = 3999.6667
$ws.Range("K102").Value = 1676.6923
$ws.Range("L102").Value = 3999.6667
$ws.Range("M102").Value = -54.69229999999993
$ws.Range("N102").Value = -7243.6667

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2508
$ws.Range("I40").Value = 2444.6785
$ws.Range("K40").Value = 2444.6785
$ws.Range("M40").Value = -2308.6785
$ws.Range("H51").Value = 43000
$ws.Range("J51").Value = 43000
$ws.Range("L51").Value = 43000
$ws.Range("N51").Value = -43956
$ws.Range("H103").Value = 14798.429
$ws.Range("J103").Value = 14798.429
$ws.Range("L103").Value = 14798.429
$ws.Range("N103").Value = -17142.429
$ws.Range("H106").Value = 12612.1
$ws.Range("J106").Value = 12612.1
$ws.Range("L106").Value = 12612.1
$ws.Range("N106").Value = -15136.1
$ws.Range("H122").Value = 6262.6875
$ws.Range("I122").Value = 3836.6365
$ws.Range("K122").Value = 11509.9095
$ws.Range("M122").Value = -9059.9095
$ws.Range("H132").Value = 31255720
$ws.Range("I132").Value = 45458150
$ws.Range("J132").Value = 10368.2
$ws.Range("K132").Value = 136374450
$ws.Range("L132").Value = 31104.6
$ws.Range("M132").Value = -136371920
$ws.Range("N132").Value = -36164.60000000001
$ws.Range("H133").Value = 78412.5
$ws.Range("J133").Value = 78412.5
$ws.Range("L133").Value = 78412.5
$ws.Range("N133").Value = -83472.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 15001750
$ws.Range("I11").Value = 15001750
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 15001750
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -15001608
$ws.Range("N11").ClearContents()
$ws.Range("H19").Value = 4666.6665
$ws.Range("J19").Value = 4666.6665
$ws.Range("L19").Value = 4666.6665
$ws.Range("N19").Value = -5014.6665
$ws.Range("H95").Value = 36333
$ws.Range("J95").Value = 36333
$ws.Range("L95").Value = 36333
$ws.Range("N95").Value = -41825
$ws.Range("H101").Value = 24483.334
$ws.Range("J101").Value = 24483.334
$ws.Range("L101").Value = 24483.334
$ws.Range("N101").Value = -30973.334
$ws.Range("H113").Value = 4985.625
$ws.Range("I113").Value = 2707.889
$ws.Range("J113").Value = 7914.143
$ws.Range("K113").Value = 8123.667
$ws.Range("L113").Value = 23742.429
$ws.Range("M113").Value = -5953.667
$ws.Range("N113").Value = -28082.429
$ws.Range("H122").Value = 2002.05
$ws.Range("I122").Value = 1646.2142
$ws.Range("J122").Value = 2832.3333
$ws.Range("K122").Value = 4938.642599999999
$ws.Range("L122").Value = 8496.999899999999
$ws.Range("M122").Value = -2488.642599999999
$ws.Range("N122").Value = -13396.9999
$ws.Range("H132").Value = 3227.5454
$ws.Range("I132").Value = 2296.8704
$ws.Range("K132").Value = 6890.611199999999
$ws.Range("M132").Value = -4360.611199999999
